# Generate Report for Handback
#
# Renames the two tracked handoff files (swapping their GUID-based names)
# and refreshes the associated handback timestamps across all three
# worksheets (Overview, zh-cn, de-de), while preserving the original
# hyperlink targets (the Address/Target URLs in the *.rels parts are not
# supposed to change - only the cell text / hyperlink display text is
# updated to reflect the new file names).

$wb = $excel.ActiveWorkbook

$oldMdA = "50527144-a073-47db-9c7c-0e38a0676b0f.md"
$oldMdB = "a45f4c37-42f6-490b-8d01-a84c223ce2ca.md"
$newMdA = "437cfbd5-767f-4178-a01b-f91116985aef.md"
$newMdB = "ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md"

$newZhCnXlf = "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf"
$newDeDeXlf = "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf"

$newZhCnHandoffTime  = "2016-03-24 10:18:15"
$newZhCnHandbackTime = "2016-03-24 10:18:40"
$newDeDeHandoffTime  = "2016-03-24 10:18:19"
$newDeDeHandbackTime = "2016-03-24 10:18:47"

function Apply-SheetUpdate($sheetName, $cellMap, $hyperlinkMap) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Plain (non-hyperlinked) cell text updates.
    foreach ($ref in $cellMap.Keys) {
        $ws.Range($ref).Value = $cellMap[$ref]
    }

    # The engine's hyperlink write-back only supports clearing the whole
    # sheet-level collection and re-adding entries (editing/deleting a
    # single pre-existing hyperlink in place is not reliable), so drop
    # every hyperlink on the sheet once and rebuild all of them with the
    # (unchanged) original target address but the refreshed display text.
    if ($hyperlinkMap.Count -gt 0) {
        $firstRef = @($hyperlinkMap.Keys)[0]
        $ws.Range($firstRef).Hyperlinks.Delete()

        foreach ($ref in $hyperlinkMap.Keys) {
            $info = $hyperlinkMap[$ref]
            $ws.Hyperlinks.Add(
                $ws.Range($ref),
                $info.Address,
                [System.Reflection.Missing]::Value,
                [System.Reflection.Missing]::Value,
                $info.Display
            )
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overviewCells = @{}
$overviewLinks = [ordered]@{
    "A2" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdA"; Display = $newMdA }
    "A3" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdB"; Display = $newMdB }
}
Apply-SheetUpdate "Overview" $overviewCells $overviewLinks

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcnCells = @{
    "E2" = $newZhCnHandoffTime
    "E3" = $newZhCnHandoffTime
    "H2" = $newZhCnHandbackTime
    "H3" = $newZhCnHandbackTime
}
$zhcnLinks = [ordered]@{
    "A2" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdA"; Display = $newMdA }
    "D2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1b2332d28c71cab1063bad17201de843e8e685a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527144-a073-47db-9c7c-0e38a0676b0f.36be47f83ad01e34c8c398658f5570c6705c86b2.zh-cn.xlf"; Display = $newZhCnXlf }
    "F2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/97e0f7bb457204b8ae7342eb169c3a1ef6437c1c/e2e/$oldMdA"; Display = $newMdA }
    "G2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ac4cfd42b31f260def9e6bf55dbf0098bb7808fb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527144-a073-47db-9c7c-0e38a0676b0f.36be47f83ad01e34c8c398658f5570c6705c86b2.zh-cn.xlf"; Display = $newZhCnXlf }
    "A3" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdB"; Display = $newMdB }
    "D3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1b2332d28c71cab1063bad17201de843e8e685a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a45f4c37-42f6-490b-8d01-a84c223ce2ca.4b9fa714f6f89c62d81c1bbc48be0a4081a2e3c1.zh-cn.xlf"; Display = $newZhCnXlf }
    "F3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/97e0f7bb457204b8ae7342eb169c3a1ef6437c1c/e2e/$oldMdB"; Display = $newMdB }
    "G3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ac4cfd42b31f260def9e6bf55dbf0098bb7808fb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a45f4c37-42f6-490b-8d01-a84c223ce2ca.4b9fa714f6f89c62d81c1bbc48be0a4081a2e3c1.zh-cn.xlf"; Display = $newZhCnXlf }
}
Apply-SheetUpdate "zh-cn" $zhcnCells $zhcnLinks

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dedeCells = @{
    "E2" = $newDeDeHandoffTime
    "E3" = $newDeDeHandoffTime
    "H2" = $newDeDeHandbackTime
    "H3" = $newDeDeHandbackTime
}
$dedeLinks = [ordered]@{
    "A2" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdA"; Display = $newMdA }
    "D2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ccf8334df655560a5a02062db52de20bac5ab217/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527144-a073-47db-9c7c-0e38a0676b0f.36be47f83ad01e34c8c398658f5570c6705c86b2.de-de.xlf"; Display = $newDeDeXlf }
    "F2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/edfa42196ece24cba48d799c4d3530edea328d6f/e2e/$oldMdA"; Display = $newMdA }
    "G2" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ae71cd1d9b0f37406ee256d3ee0b6773708244c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527144-a073-47db-9c7c-0e38a0676b0f.36be47f83ad01e34c8c398658f5570c6705c86b2.de-de.xlf"; Display = $newDeDeXlf }
    "A3" = @{ Address = "https://github.com/OpenLocalizationTest/oltest/blob/7c7e99c95fcfbb66d19cc97419ef4d5077845f78/e2e/$oldMdB"; Display = $newMdB }
    "D3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ccf8334df655560a5a02062db52de20bac5ab217/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a45f4c37-42f6-490b-8d01-a84c223ce2ca.4b9fa714f6f89c62d81c1bbc48be0a4081a2e3c1.de-de.xlf"; Display = $newDeDeXlf }
    "F3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/edfa42196ece24cba48d799c4d3530edea328d6f/e2e/$oldMdB"; Display = $newMdB }
    "G3" = @{ Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ae71cd1d9b0f37406ee256d3ee0b6773708244c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a45f4c37-42f6-490b-8d01-a84c223ce2ca.4b9fa714f6f89c62d81c1bbc48be0a4081a2e3c1.de-de.xlf"; Display = $newDeDeXlf }
}
Apply-SheetUpdate "de-de" $dedeCells $dedeLinks
